# Removed Extension Payments Tax Type from execution.
# This script replays a Katalon test-suite re-run: every row's "Date" stamp
# is refreshed, and every row whose PaymentType is "Extension Payments" is
# flagged as DONOTRUN/DoNotRun (Execute column) with Result flipped to Fail.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Estimated
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Thu Dec 07 22:11:28 EST 2023"
$ws.Range("B3").Value = "Thu Dec 07 22:12:21 EST 2023"
$ws.Range("B4").Value = "Thu Dec 07 22:13:15 EST 2023"
$ws.Range("B5").Value = "Thu Dec 07 22:14:09 EST 2023"
$ws.Range("B6").Value = "Thu Dec 07 22:15:03 EST 2023"
$ws.Range("B7").Value = "Thu Dec 07 22:15:56 EST 2023"

# ---------------------------------------------------------------------------
# Existing  (loses tabSelected -- handled by not activating it last)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Thu Dec 07 22:16:50 EST 2023"
$ws.Range("B3").Value = "Thu Dec 07 22:17:43 EST 2023"
$ws.Range("B4").Value = "Thu Dec 07 22:18:36 EST 2023"
$ws.Range("B5").Value = "Thu Dec 07 22:19:29 EST 2023"
$ws.Range("B6").Value = "Thu Dec 07 22:20:21 EST 2023"
$ws.Range("B7").Value = "Thu Dec 07 22:21:14 EST 2023"
$ws.Range("B8").Value = "Thu Dec 07 22:22:08 EST 2023"
$ws.Range("B9").Value = "Thu Dec 07 22:23:02 EST 2023"
$ws.Range("B10").Value = "Thu Dec 07 22:23:56 EST 2023"
$ws.Range("B11").Value = "Thu Dec 07 22:24:51 EST 2023"
$ws.Range("B12").Value = "Thu Dec 07 22:25:45 EST 2023"

# ---------------------------------------------------------------------------
# Extension  -- every row pays "Extension Payments" -> disable them all
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Extension")
$ws.Activate()

$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Oct 31 16:41:50 EDT 2023"
$ws.Range("C2").Value = "DONOTRUN"

$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Oct 31 16:42:47 EDT 2023"
$ws.Range("C3").Value = "DONOTRUN"

$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Tue Oct 31 16:43:40 EDT 2023"
$ws.Range("C4").Value = "DONOTRUN"

$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Tue Oct 31 16:44:33 EDT 2023"
$ws.Range("C5").Value = "DONOTRUN"

$ws.Range("A6").Value = "Fail"
$ws.Range("B6").Value = "Tue Oct 31 16:45:26 EDT 2023"
$ws.Range("C6").Value = "DONOTRUN"

$ws.Range("A7").Value = "Fail"
$ws.Range("B7").Value = "Tue Oct 31 16:46:19 EDT 2023"
$ws.Range("C7").Value = "DONOTRUN"

$ws.Columns.Item(3).ColumnWidth = 18.25
$ws.Range("C2:C7").Select()

# ---------------------------------------------------------------------------
# NewTaxReturn
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Thu Dec 07 22:26:39 EST 2023"
$ws.Range("B3").Value = "Thu Dec 07 22:27:33 EST 2023"
$ws.Range("B4").Value = "Thu Dec 07 22:28:26 EST 2023"
$ws.Range("B5").Value = "Thu Dec 07 22:29:19 EST 2023"
$ws.Range("B6").Value = "Thu Dec 07 22:30:13 EST 2023"
$ws.Range("B7").Value = "Thu Dec 07 22:31:07 EST 2023"
$ws.Range("B8").Value = "Thu Dec 07 22:32:01 EST 2023"
$ws.Range("B9").Value = "Thu Dec 07 22:32:54 EST 2023"
$ws.Range("B10").Value = "Thu Dec 07 22:33:47 EST 2023"
$ws.Range("B11").Value = "Thu Dec 07 22:34:40 EST 2023"
$ws.Range("B12").Value = "Thu Dec 07 22:35:32 EST 2023"
$ws.Range("B13").Value = "Thu Dec 07 22:36:25 EST 2023"
$ws.Range("B14").Value = "Thu Dec 07 22:37:19 EST 2023"
$ws.Range("B15").Value = "Thu Dec 07 22:38:12 EST 2023"
$ws.Range("B16").Value = "Thu Dec 07 22:39:05 EST 2023"

# ---------------------------------------------------------------------------
# Personal -- row 4 (Extension Payments) disabled
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Personal")
$ws.Activate()
$ws.Range("C4").Value = "DoNotRun"
$ws.Columns.Item(3).ColumnWidth = 25.42
$ws.Range("C4").Select()

# ---------------------------------------------------------------------------
# Personal_IND -- row 3 (Extension Payments) disabled; becomes active sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Activate()

$ws.Range("B2").Value = "Fri Dec 08 10:50:48 EST 2023"

$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Thu Dec 07 22:41:41 EST 2023"
$ws.Range("C3").Value = "DoNotRun"

$ws.Range("B4").Value = "Fri Dec 08 10:51:38 EST 2023"
$ws.Range("B5").Value = "Fri Dec 08 10:52:28 EST 2023"
$ws.Range("B6").Value = "Fri Dec 08 10:53:18 EST 2023"

$ws.Columns.Item(3).ColumnWidth = 14.25
$ws.Range("C3").Select()

# ---------------------------------------------------------------------------
# Personal_JNT -- row 3 (Extension Payments) disabled
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Activate()

$ws.Range("B2").Value = "Fri Dec 08 10:54:09 EST 2023"

$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Thu Dec 07 22:46:40 EST 2023"
$ws.Range("C3").Value = "DoNotRun"

$ws.Range("B4").Value = "Fri Dec 08 10:55:08 EST 2023"
$ws.Range("B5").Value = "Fri Dec 08 10:56:04 EST 2023"
$ws.Range("B6").Value = "Fri Dec 08 10:57:02 EST 2023"

$ws.Columns.Item(3).ColumnWidth = 15.42
$ws.Range("C3").Select()

# ---------------------------------------------------------------------------
# Personal_EL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Thu Dec 07 22:39:59 EST 2023"

# ---------------------------------------------------------------------------
# Final active sheet / selection: Personal_IND (tabSelected, activeTab)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Activate()
$ws.Range("C3").Select()
